$wb = $excel.ActiveWorkbook

# --- 1. Add a new row to the "Test Cases" sheet for the new GoogleSearch1Test case ---
$testCases = $wb.Worksheets.Item("Test Cases")
$testCases.Cells.Item(3, 1).Value = "GoogleSearch1Test"
$testCases.Cells.Item(3, 2).Value = "xxxxxxxxx"
$testCases.Cells.Item(3, 3).Value = "Y"
$testCases.Cells.Item(3, 4).Value = "PASS"
$testCases.Cells.Item(3, 4).WrapText = $true
$testCases.Range("A3").Select()

# --- 2. Duplicate the "GoogleSearchTest" sheet into a new "GoogleSearch1Test" sheet ---
$googleSearchTest = $wb.Worksheets.Item("GoogleSearchTest")
$googleSearchTest.Copy($null, $googleSearchTest)
$newSheet = $wb.Worksheets.Item($googleSearchTest.Index + 1)
$newSheet.Name = "GoogleSearch1Test"

# Re-apply the wrap formatting on the new sheet's F2 cell (Results/PASS for Data4=Google)
$newSheet.Cells.Item(2, 6).WrapText = $true

# --- 3. Update selections: whole-table range selected on both GoogleSearchTest sheets ---
$googleSearchTest.Range("A1:F3").Select()
$googleSearchTest.Range("F3").Activate()

$newSheet.Range("A1:F3").Select()
$newSheet.Activate()
